$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly record was added to the "Apio" price series. It belongs
# chronologically at row 282 (pushing the existing rows 282-346 down to
# 283-347), so insert a fresh row there and fill it in with the new data.
$ws.Rows.Item(282).Insert()

$ws.Cells.Item(282, 1).Value = 10
$ws.Cells.Item(282, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(282, 3).Value = "La Araucanía"
$ws.Cells.Item(282, 4).Value = 44785
$ws.Cells.Item(282, 5).Value = 9
$ws.Cells.Item(282, 6).Value = 100112017
$ws.Cells.Item(282, 7).Value = "Apio"
$ws.Cells.Item(282, 8).Value = "Americana (o)"
$ws.Cells.Item(282, 9).Value = "Primera"
$ws.Cells.Item(282, 10).Value = 50
$ws.Cells.Item(282, 11).Value = 11000
$ws.Cells.Item(282, 12).Value = 12000
$ws.Cells.Item(282, 13).Value = 11400
$ws.Cells.Item(282, 14).Value = "`$/docena de matas"
$ws.Cells.Item(282, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(282, 16).Value = 1900
$ws.Cells.Item(282, 17).Value = 6
$ws.Cells.Item(282, 18).Value = "Hortaliza"

# Match the date-number format used by the other rows in column D.
$ws.Cells.Item(282, 4).NumberFormat = $ws.Cells.Item(283, 4).NumberFormat
